# update common template imports
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("survey")

# Row 3 - date picker template now lives in the shared assets folder
$ws.Range("F3").Value = "../../../../assets/templates/custom_date_picker.handlebars"

# Row 1 - new "comment" header in column G
$ws.Range("G1").Value = "comment"

# Row 3 - comment explaining the shared-asset template import
$ws.Range("G3").Value = "reusable templates are automatically copied to asset folder"

# Row 4 - datetime picker template is imported directly from the same folder, plus a comment
$ws.Range("F4").Value = "custom_datetime_picker.handlebars"
$ws.Range("G4").Value = "templates in same folder can be directly imported"

# Row heights adjust to fit the new comment text
$ws.Rows.Item(3).RowHeight = 63.75
$ws.Rows.Item(4).RowHeight = 25.5

# New column G is widened to show the comment text
$ws.Columns.Item(7).ColumnWidth = 26 - 5/6

# Selection moves to the newly added cell
$ws.Range("G5").Select() | Out-Null
